# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 79 (pushing the existing rows
# 79:153 down to 80:154) in the "Arándano (blue)" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 79; Excel copies the row-above's
# formatting (incl. the date style on column D) the same way a manual
# "Insert Sheet Rows" does.
$ws.Rows.Item(79).Insert()

$newRow = 79

$ws.Cells.Item($newRow, 1).Value = 10
$ws.Cells.Item($newRow, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($newRow, 3).Value = "La Araucanía"
$ws.Cells.Item($newRow, 4).Value = (Get-Date -Year 2023 -Month 12 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item($newRow, 5).Value = 9
$ws.Cells.Item($newRow, 6).Value = "Fruta"
$ws.Cells.Item($newRow, 7).Value = 100101
$ws.Cells.Item($newRow, 8).Value = "Berries"
$ws.Cells.Item($newRow, 9).Value = 100101001
$ws.Cells.Item($newRow, 10).Value = "Arándano (blue)"
$ws.Cells.Item($newRow, 11).Value = "Sin especificar"
$ws.Cells.Item($newRow, 12).Value = "Primera"
$ws.Cells.Item($newRow, 13).Value = 800
$ws.Cells.Item($newRow, 14).Value = 3300
$ws.Cells.Item($newRow, 15).Value = 3400
$ws.Cells.Item($newRow, 16).Value = 3338
$ws.Cells.Item($newRow, 17).Value = "`$/kilo"
$ws.Cells.Item($newRow, 18).Value = "Región del Maule"
$ws.Cells.Item($newRow, 19).Value = 3338
$ws.Cells.Item($newRow, 20).Value = 1
